$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F ("dSF")
$updates = @{
    6  = 5
    7  = -3
    12 = -3
    13 = 3
    14 = 3
    17 = -1
    19 = -1
    22 = -1
    25 = 0
    29 = 2
    30 = -3
    33 = -2
    37 = 3
    38 = 2
    50 = -1
    58 = -1
    61 = 0
    70 = -9
    71 = -8
    74 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
